$d = $word.ActiveDocument

# Update the date line at the top of the document (single occurrence in
# the doc, so a plain document-wide Find/Replace is unambiguous).
$d.Content.Find.Execute("2025-06-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-23 Monday", 2) | Out-Null

# Update each division problem in the practice table. Several problem
# strings (e.g. "13÷3=") occur more than once in the table, so a global
# Find/Replace would be ambiguous. Instead target each table cell by its
# (row, column) position directly and overwrite just the run text in that
# cell (Start .. Start+len(oldText)), which leaves every other cell, and
# the existing run formatting (font/size) in the edited cell, untouched.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cellStart = $cell.Range.Start
$old = "77÷6="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "96÷9="

$cell = $t.Cell(1, 2)
$cellStart = $cell.Range.Start
$old = "50÷7="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "59÷4="

$cell = $t.Cell(1, 3)
$cellStart = $cell.Range.Start
$old = "11÷2="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "72÷7="

$cell = $t.Cell(1, 4)
$cellStart = $cell.Range.Start
$old = "94÷2="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "86÷8="

$cell = $t.Cell(1, 5)
$cellStart = $cell.Range.Start
$old = "54÷4="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "91÷4="

$cell = $t.Cell(5, 1)
$cellStart = $cell.Range.Start
$old = "49÷6="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "98÷2="

$cell = $t.Cell(5, 2)
$cellStart = $cell.Range.Start
$old = "14÷9="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "21÷2="

$cell = $t.Cell(5, 3)
$cellStart = $cell.Range.Start
$old = "42÷8="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "13÷8="

$cell = $t.Cell(5, 4)
$cellStart = $cell.Range.Start
$old = "41÷5="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "55÷9="

$cell = $t.Cell(5, 5)
$cellStart = $cell.Range.Start
$old = "99÷5="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "36÷8="

$cell = $t.Cell(9, 1)
$cellStart = $cell.Range.Start
$old = "93÷5="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "87÷6="

$cell = $t.Cell(9, 2)
$cellStart = $cell.Range.Start
$old = "94÷8="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "23÷8="

$cell = $t.Cell(9, 3)
$cellStart = $cell.Range.Start
$old = "40÷5="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "78÷4="

$cell = $t.Cell(9, 4)
$cellStart = $cell.Range.Start
$old = "79÷4="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "50÷6="

$cell = $t.Cell(9, 5)
$cellStart = $cell.Range.Start
$old = "13÷3="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "72÷6="

$cell = $t.Cell(13, 1)
$cellStart = $cell.Range.Start
$old = "96÷2="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "80÷2="

$cell = $t.Cell(13, 2)
$cellStart = $cell.Range.Start
$old = "33÷4="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "82÷3="

$cell = $t.Cell(13, 3)
$cellStart = $cell.Range.Start
$old = "13÷3="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "13÷6="

$cell = $t.Cell(13, 4)
$cellStart = $cell.Range.Start
$old = "17÷7="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "14÷5="

$cell = $t.Cell(13, 5)
$cellStart = $cell.Range.Start
$old = "75÷9="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "76÷9="

$cell = $t.Cell(17, 1)
$cellStart = $cell.Range.Start
$old = "64÷8="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "90÷8="

$cell = $t.Cell(17, 2)
$cellStart = $cell.Range.Start
$old = "21÷5="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "28÷7="

$cell = $t.Cell(17, 3)
$cellStart = $cell.Range.Start
$old = "71÷7="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "76÷3="

$cell = $t.Cell(17, 4)
$cellStart = $cell.Range.Start
$old = "32÷2="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "83÷8="

$cell = $t.Cell(17, 5)
$cellStart = $cell.Range.Start
$old = "51÷9="
$rng = $d.Range($cellStart, $cellStart + $old.Length)
$rng.Text = "70÷8="
